# This script applies updated odds values to Sheet1 of the workbook,
# matching the target diff (rows 3,4,5,9,11,12,13,15,16,17,18,19,25,28,31,33,35,36).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 2.9
$ws.Range("I3").Value = 2.55
$ws.Range("K3").Value = 10
$ws.Range("AG3").Value = 10
$ws.Range("AI3").Value = 21

# Row 4
$ws.Range("G4").Value = 8.5

# Row 5
$ws.Range("G5").Value = 3.25
$ws.Range("I5").Value = 2.7
$ws.Range("J5").Value = 1.14
$ws.Range("K5").Value = 5.5
$ws.Range("AG5").Value = 11

# Row 9
$ws.Range("J9").Value = 1.1
$ws.Range("K9").Value = 7

# Row 11
$ws.Range("G11").Value = 1.98
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 3.7
$ws.Range("L11").Value = 1.45
$ws.Range("M11").Value = 2.37
$ws.Range("N11").Value = 2.32
$ws.Range("O11").Value = 1.47
$ws.Range("P11").Value = 1.52
$ws.Range("Q11").Value = 2.22
$ws.Range("R11").Value = 2.07
$ws.Range("S11").Value = 1.6
$ws.Range("T11").Value = 5.5
$ws.Range("U11").Value = 8
$ws.Range("V11").Value = 9.25
$ws.Range("W11").Value = 17
$ws.Range("X11").Value = 19.5
$ws.Range("Y11").Value = 40
$ws.Range("Z11").Value = 7.1
$ws.Range("AA11").Value = 6.4
$ws.Range("AB11").Value = 20
$ws.Range("AC11").Value = 120
$ws.Range("AE11").Value = 8.25
$ws.Range("AF11").Value = 17.5
$ws.Range("AG11").Value = 13.5
$ws.Range("AH11").Value = 55
$ws.Range("AI11").Value = 45
$ws.Range("AJ11").Value = 60

# Row 12
$ws.Range("G12").Value = 1.87
$ws.Range("H12").Value = 3.05
$ws.Range("I12").Value = 4.4
$ws.Range("L12").Value = 1.53
$ws.Range("M12").Value = 2.2
$ws.Range("N12").Value = 2.52
$ws.Range("O12").Value = 1.4
$ws.Range("P12").Value = 1.57
$ws.Range("Q12").Value = 2.1
$ws.Range("R12").Value = 2.25
$ws.Range("S12").Value = 1.5
$ws.Range("T12").Value = 4.85
$ws.Range("U12").Value = 7.2
$ws.Range("V12").Value = 9.5
$ws.Range("W12").Value = 15
$ws.Range("X12").Value = 20
$ws.Range("Y12").Value = 50
$ws.Range("Z12").Value = 6
$ws.Range("AA12").Value = 6.4
$ws.Range("AB12").Value = 23
$ws.Range("AC12").Value = 175
$ws.Range("AE12").Value = 8.75
$ws.Range("AF12").Value = 22
$ws.Range("AG12").Value = 16
$ws.Range("AH12").Value = 80
$ws.Range("AI12").Value = 60
$ws.Range("AJ12").Value = 80

# Row 13
$ws.Range("G13").Value = 2.15
$ws.Range("I13").Value = 3.4
$ws.Range("W13").Value = 19
$ws.Range("AA13").Value = 6.5
$ws.Range("AD13").Value = 401
$ws.Range("AF13").Value = 17

# Row 15
$ws.Range("G15").Value = 5.3
$ws.Range("H15").Value = 3.65
$ws.Range("I15").Value = 1.55
$ws.Range("N15").Value = 1.85
$ws.Range("O15").Value = 1.75
$ws.Range("P15").Value = 1.38
$ws.Range("Q15").Value = 2.47
$ws.Range("T15").Value = 11
$ws.Range("U15").Value = 25
$ws.Range("V15").Value = 14
$ws.Range("W15").Value = 75
$ws.Range("X15").Value = 45
$ws.Range("Y15").Value = 45
$ws.Range("Z15").Value = 9.25
$ws.Range("AA15").Value = 6.3
$ws.Range("AB15").Value = 14.5
$ws.Range("AC15").Value = 65
$ws.Range("AD15").Value = 500
$ws.Range("AE15").Value = 5.4
$ws.Range("AF15").Value = 5.9
$ws.Range("AG15").Value = 6.9
$ws.Range("AH15").Value = 9.25
$ws.Range("AI15").Value = 10.75
$ws.Range("AJ15").Value = 22

# Row 16
$ws.Range("G16").Value = 2.18
$ws.Range("H16").Value = 3.45
$ws.Range("I16").Value = 2.77
$ws.Range("N16").Value = 1.75
$ws.Range("O16").Value = 1.87
$ws.Range("P16").Value = 1.38
$ws.Range("Q16").Value = 2.47
$ws.Range("T16").Value = 7.3
$ws.Range("U16").Value = 9.25
$ws.Range("V16").Value = 7.7
$ws.Range("W16").Value = 17
$ws.Range("X16").Value = 14
$ws.Range("Y16").Value = 20
$ws.Range("Z16").Value = 11.25
$ws.Range("AA16").Value = 6
$ws.Range("AB16").Value = 11.25
$ws.Range("AC16").Value = 40
$ws.Range("AD16").Value = 250
$ws.Range("AE16").Value = 8.25
$ws.Range("AF16").Value = 12
$ws.Range("AG16").Value = 8.75
$ws.Range("AH16").Value = 25
$ws.Range("AI16").Value = 18
$ws.Range("AJ16").Value = 23

# Row 17
$ws.Range("J17").Value = 1.05
$ws.Range("K17").Value = 11
$ws.Range("N17").Value = 1.83
$ws.Range("O17").Value = 1.98
$ws.Range("P17").Value = 1.36
$ws.Range("Q17").Value = 3
$ws.Range("Z17").Value = 11

# Row 18
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 2.9
$ws.Range("J18").Value = 1.08
$ws.Range("K18").Value = 8
$ws.Range("N18").Value = 2.25
$ws.Range("O18").Value = 1.62
$ws.Range("R18").Value = 1.95
$ws.Range("S18").Value = 1.8
$ws.Range("T18").Value = 7.5
$ws.Range("X18").Value = 23
$ws.Range("Z18").Value = 8
$ws.Range("AD18").Value = 301
$ws.Range("AE18").Value = 8
$ws.Range("AI18").Value = 26
$ws.Range("AJ18").Value = 41

# Row 19
$ws.Range("I19").Value = 3.4
$ws.Range("K19").Value = 12
$ws.Range("Z19").Value = 12
$ws.Range("AH19").Value = 41

# Row 25
$ws.Range("J25").Value = 1.02
$ws.Range("K25").Value = 11

# Row 28
$ws.Range("G28").Value = 6.5
$ws.Range("H28").Value = 6
$ws.Range("I28").Value = 1.27
$ws.Range("J28").Value = 1.01
$ws.Range("K28").Value = 34
$ws.Range("L28").Value = 1.05
$ws.Range("M28").Value = 11
$ws.Range("N28").Value = 1.2
$ws.Range("O28").Value = 4.33
$ws.Range("P28").Value = 1.14
$ws.Range("Q28").Value = 5.5
$ws.Range("R28").Value = 1.44
$ws.Range("S28").Value = 2.63
$ws.Range("T28").Value = 34
$ws.Range("U28").Value = 51
$ws.Range("V28").Value = 23
$ws.Range("W28").Value = 81
$ws.Range("X28").Value = 41
$ws.Range("Y28").Value = 34
$ws.Range("Z28").Value = 34
$ws.Range("AA28").Value = 15
$ws.Range("AB28").Value = 17
$ws.Range("AC28").Value = 34
$ws.Range("AD28").Value = 81
$ws.Range("AE28").Value = 17
$ws.Range("AF28").Value = 12
$ws.Range("AG28").Value = 10
$ws.Range("AH28").Value = 12
$ws.Range("AI28").Value = 10
$ws.Range("AJ28").Value = 17

# Row 31
$ws.Range("K31").Value = 13

# Row 33
$ws.Range("G33").Value = 6.5
$ws.Range("H33").Value = 4.2
$ws.Range("I33").Value = 1.48
$ws.Range("L33").Value = 1.29
$ws.Range("M33").Value = 3.5
$ws.Range("R33").Value = 2
$ws.Range("S33").Value = 1.75
$ws.Range("AC33").Value = 67
$ws.Range("AD33").Value = 401

# Row 35
$ws.Range("H35").Value = 4
$ws.Range("I35").Value = 3.8
$ws.Range("R35").Value = 1.44
$ws.Range("S35").Value = 2.63
$ws.Range("AA35").Value = 8

# Row 36
$ws.Range("G36").Value = 2.7
$ws.Range("I36").Value = 2.63
$ws.Range("T36").Value = 8
$ws.Range("U36").Value = 13
$ws.Range("V36").Value = 11
$ws.Range("W36").Value = 26
$ws.Range("X36").Value = 23
$ws.Range("AE36").Value = 8
$ws.Range("AF36").Value = 12
$ws.Range("AH36").Value = 26
